# Apply the edits described by the diff to the workbook.
$wb = $excel.ActiveWorkbook

# ---- Metadata sheet updates ----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-allowed-amount-primary-physician"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet updates ----
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element; its Constraint(s) cell (AI) is cleared.
$elements.Range("AI2").Value = ""

# Row 5 is "Extension.url"; its Fixed Value (Q) mirrors the canonical extension URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-allowed-amount-primary-physician"
